# Remove the paragraph-mark run-properties font hint
# (<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr> living directly under
# <w:pPr>) from every paragraph that has one. This does NOT touch the
# run-level <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr> markup,
# which must stay exactly as-is.

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $xml = $p.Range.WordOpenXML

    if ($xml -match '<w:body>(<w:p\b.*?</w:p>|<w:p\b[^>]*/>)') {
        $frag = $matches[1]

        # Only touch paragraphs whose pPr directly carries the
        # eastAsia-hint rFonts run-properties block.
        if ($frag -match '<w:pPr>((?:(?!</w:pPr>).)*?)<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>((?:(?!</w:pPr>).)*?)</w:pPr>') {

            # Strip the synthetic paragraph/text ids the WordOpenXML
            # getter stamps on (they aren't present in the source and
            # shouldn't be introduced).
            $newFrag = $frag -replace ' w14:paraId="[0-9A-Fa-f]+"', ''
            $newFrag = $newFrag -replace ' w14:textId="[0-9A-Fa-f]+"', ''

            # Drop just the paragraph-mark rPr; if pPr becomes empty,
            # drop pPr entirely too.
            $newFrag = $newFrag -replace '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>', ''
            $newFrag = $newFrag -replace '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>', '</w:pPr>'

            [void]$p.Range.InsertXML($newFrag)
        }
    }
}

Write-Output "done"
